# Automatische test-sync: 2025-08-14 20:39:50
#
# Adds a new "Inkoop / Bestellingen" mail-log entry (row 11) to the Logs
# sheet, updates the conditional formatting ranges to include it, adds the
# matching aggregate row (row 4) to the Dashboard sheet, and extends the
# bar chart's category/value series references to cover the new row.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new row -----------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A11").Value = "Nieuwe bestelling"
$logs.Range("B11").Value = "planning@testbedrijf123.nl"
$logs.Range("C11").Value = "Wil je 200 stuks M8-bouten bestellen bij onze leverancier?"
$logs.Range("D11").Value = "Inkoop / Bestellingen"
$logs.Range("E11").Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@testbedrijf123.nl."
$logs.Range("F11").Value = "2025-08-14 20:39:34"
$logs.Range("G11").Value = "Nee"
$logs.Range("H11").Value = "Ja"
$logs.Range("I11").Value = "Nee"
$logs.Range("J11").Value = "Nee"

# --- Logs sheet: extend conditional formatting ranges to row 11 -------
$cfColumns = @("D", "G", "H", "I", "J")
foreach ($col in $cfColumns) {
    $oldRange = $logs.Range($col + "2:" + $col + "10")
    $newRange = $logs.Range($col + "2:" + $col + "11")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard sheet: append aggregate row for the new category -------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A4").Value = "Inkoop / Bestellingen"
$dashboard.Range("B4").Value = 1

# --- Dashboard chart: extend category/value series to row 4 -----------
$chartObj = $dashboard.ChartObjects().Item(1)
$series = $chartObj.Chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$4,'Dashboard'!`$B`$2:`$B`$4,1)"
